# Refresh simulated game-outcome transition probabilities for Saint Louis_A
# after adding more simulated games and speeding up the simulate-game logic.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = New-Object 'object[,]' 18,18
$data[0,0] = 0.176056338028169
$data[0,1] = 0.5950704225352113
$data[0,2] = 0
$data[0,3] = 0
$data[0,4] = 0
$data[0,5] = 0
$data[0,6] = 0
$data[0,7] = 0
$data[0,8] = 0.01408450704225352
$data[0,9] = 0
$data[0,10] = 0
$data[0,11] = 0
$data[0,12] = 0
$data[0,13] = 0
$data[0,14] = 0.147887323943662
$data[0,15] = 0
$data[0,16] = 0
$data[0,17] = 0.06690140845070422
$data[1,0] = 0.01176470588235294
$data[1,1] = 0.01176470588235294
$data[1,2] = 0
$data[1,3] = 0
$data[1,4] = 0
$data[1,5] = 0
$data[1,6] = 0
$data[1,7] = 0
$data[1,8] = 0.02352941176470588
$data[1,9] = 0
$data[1,10] = 0
$data[1,11] = 0
$data[1,12] = 0
$data[1,13] = 0
$data[1,14] = 0.7529411764705882
$data[1,15] = 0
$data[1,16] = 0
$data[1,17] = 0.2
$data[2,0] = 0
$data[2,1] = 0
$data[2,2] = 0
$data[2,3] = 0
$data[2,4] = 0
$data[2,5] = 0
$data[2,6] = 0
$data[2,7] = 0
$data[2,8] = 0
$data[2,9] = 0
$data[2,10] = 0
$data[2,11] = 0
$data[2,12] = 0
$data[2,13] = 0
$data[2,14] = 0.7291666666666666
$data[2,15] = 0
$data[2,16] = 0
$data[2,17] = 0.2708333333333333
$data[3,0] = 0
$data[3,1] = 0
$data[3,2] = 0
$data[3,3] = 0
$data[3,4] = 0
$data[3,5] = 0
$data[3,6] = 0
$data[3,7] = 0
$data[3,8] = 0
$data[3,9] = 0
$data[3,10] = 0
$data[3,11] = 0
$data[3,12] = 0
$data[3,13] = 0
$data[3,14] = 0
$data[3,15] = 0
$data[3,16] = 0
$data[3,17] = 0
$data[4,0] = 0.06382978723404255
$data[4,1] = 0
$data[4,2] = 0.01595744680851064
$data[4,3] = 0
$data[4,4] = 0.03723404255319149
$data[4,5] = 0
$data[4,6] = 0
$data[4,7] = 0
$data[4,8] = 0.2446808510638298
$data[4,9] = 0
$data[4,10] = 0
$data[4,11] = 0
$data[4,12] = 0
$data[4,13] = 0.02127659574468085
$data[4,14] = 0
$data[4,15] = 0.1436170212765958
$data[4,16] = 0.09042553191489362
$data[4,17] = 0.3829787234042553
$data[5,0] = 0.1185567010309278
$data[5,1] = 0
$data[5,2] = 0.02061855670103093
$data[5,3] = 0
$data[5,4] = 0.04639175257731959
$data[5,5] = 0
$data[5,6] = 0
$data[5,7] = 0
$data[5,8] = 0.1134020618556701
$data[5,9] = 0
$data[5,10] = 0
$data[5,11] = 0
$data[5,12] = 0
$data[5,13] = 0.03092783505154639
$data[5,14] = 0
$data[5,15] = 0.1494845360824742
$data[5,16] = 0.07216494845360824
$data[5,17] = 0.4484536082474227
$data[6,0] = 0.1296296296296296
$data[6,1] = 0
$data[6,2] = 0.01851851851851852
$data[6,3] = 0
$data[6,4] = 0.06084656084656084
$data[6,5] = 0
$data[6,6] = 0
$data[6,7] = 0
$data[6,8] = 0.1164021164021164
$data[6,9] = 0
$data[6,10] = 0
$data[6,11] = 0
$data[6,12] = 0
$data[6,13] = 0.01587301587301587
$data[6,14] = 0
$data[6,15] = 0.1402116402116402
$data[6,16] = 0.09259259259259259
$data[6,17] = 0.4259259259259259
$data[7,0] = 0.1512195121951219
$data[7,1] = 0
$data[7,2] = 0.01951219512195122
$data[7,3] = 0
$data[7,4] = 0.03902439024390244
$data[7,5] = 0
$data[7,6] = 0
$data[7,7] = 0
$data[7,8] = 0.07804878048780488
$data[7,9] = 0
$data[7,10] = 0
$data[7,11] = 0
$data[7,12] = 0
$data[7,13] = 0.01951219512195122
$data[7,14] = 0
$data[7,15] = 0.1804878048780488
$data[7,16] = 0.08780487804878048
$data[7,17] = 0.424390243902439
$data[8,0] = 0.108133971291866
$data[8,1] = 0
$data[8,2] = 0.03062200956937799
$data[8,3] = 0
$data[8,4] = 0.07942583732057416
$data[8,5] = 0
$data[8,6] = 0
$data[8,7] = 0
$data[8,8] = 0.09856459330143541
$data[8,9] = 0
$data[8,10] = 0
$data[8,11] = 0
$data[8,12] = 0
$data[8,13] = 0.02296650717703349
$data[8,14] = 0
$data[8,15] = 0.215311004784689
$data[8,16] = 0.08899521531100478
$data[8,17] = 0.3559808612440191
$data[9,0] = 0
$data[9,1] = 0
$data[9,2] = 0
$data[9,3] = 0
$data[9,4] = 0
$data[9,5] = 0.1423728813559322
$data[9,6] = 0
$data[9,7] = 0
$data[9,8] = 0.0847457627118644
$data[9,9] = 0.1864406779661017
$data[9,10] = 0.576271186440678
$data[9,11] = 0
$data[9,12] = 0
$data[9,13] = 0
$data[9,14] = 0
$data[9,15] = 0
$data[9,16] = 0
$data[9,17] = 0.01016949152542373
$data[10,0] = 0
$data[10,1] = 0
$data[10,2] = 0
$data[10,3] = 0
$data[10,4] = 0
$data[10,5] = 0.7613636363636364
$data[10,6] = 0
$data[10,7] = 0
$data[10,8] = 0.1761363636363636
$data[10,9] = 0
$data[10,10] = 0.02840909090909091
$data[10,11] = 0
$data[10,12] = 0
$data[10,13] = 0
$data[10,14] = 0
$data[10,15] = 0
$data[10,16] = 0
$data[10,17] = 0.03409090909090909
$data[11,0] = 0
$data[11,1] = 0
$data[11,2] = 0
$data[11,3] = 0
$data[11,4] = 0
$data[11,5] = 0.7027027027027027
$data[11,6] = 0
$data[11,7] = 0
$data[11,8] = 0.1891891891891892
$data[11,9] = 0
$data[11,10] = 0
$data[11,11] = 0
$data[11,12] = 0
$data[11,13] = 0
$data[11,14] = 0
$data[11,15] = 0
$data[11,16] = 0
$data[11,17] = 0.1081081081081081
$data[12,0] = 0
$data[12,1] = 0
$data[12,2] = 0
$data[12,3] = 0
$data[12,4] = 0
$data[12,5] = 1
$data[12,6] = 0
$data[12,7] = 0
$data[12,8] = 0
$data[12,9] = 0
$data[12,10] = 0
$data[12,11] = 0
$data[12,12] = 0
$data[12,13] = 0
$data[12,14] = 0
$data[12,15] = 0
$data[12,16] = 0
$data[12,17] = 0
$data[13,0] = 0
$data[13,1] = 0
$data[13,2] = 0
$data[13,3] = 0
$data[13,4] = 0.0155440414507772
$data[13,5] = 0
$data[13,6] = 0.1450777202072539
$data[13,7] = 0.05699481865284974
$data[13,8] = 0.3212435233160622
$data[13,9] = 0.06217616580310881
$data[13,10] = 0
$data[13,11] = 0.02590673575129534
$data[13,12] = 0
$data[13,13] = 0.05699481865284974
$data[13,14] = 0
$data[13,15] = 0
$data[13,16] = 0
$data[13,17] = 0.3160621761658031
$data[14,0] = 0
$data[14,1] = 0
$data[14,2] = 0
$data[14,3] = 0
$data[14,4] = 0.0303030303030303
$data[14,5] = 0
$data[14,6] = 0.196969696969697
$data[14,7] = 0.1262626262626263
$data[14,8] = 0.3181818181818182
$data[14,9] = 0.1313131313131313
$data[14,10] = 0
$data[14,11] = 0.0202020202020202
$data[14,12] = 0
$data[14,13] = 0.05555555555555555
$data[14,14] = 0
$data[14,15] = 0
$data[14,16] = 0
$data[14,17] = 0.1212121212121212
$data[15,0] = 0
$data[15,1] = 0
$data[15,2] = 0
$data[15,3] = 0
$data[15,4] = 0.01626016260162602
$data[15,5] = 0
$data[15,6] = 0.1815718157181572
$data[15,7] = 0.1002710027100271
$data[15,8] = 0.3983739837398374
$data[15,9] = 0.1002710027100271
$data[15,10] = 0
$data[15,11] = 0.02439024390243903
$data[15,12] = 0.002710027100271003
$data[15,13] = 0.06233062330623306
$data[15,14] = 0
$data[15,15] = 0
$data[15,16] = 0
$data[15,17] = 0.1138211382113821
$data[16,0] = 0
$data[16,1] = 0
$data[16,2] = 0
$data[16,3] = 0
$data[16,4] = 0.01136363636363636
$data[16,5] = 0
$data[16,6] = 0.1534090909090909
$data[16,7] = 0.1079545454545455
$data[16,8] = 0.4034090909090909
$data[16,9] = 0.09090909090909091
$data[16,10] = 0
$data[16,11] = 0.01136363636363636
$data[16,12] = 0.005681818181818182
$data[16,13] = 0.09090909090909091
$data[16,14] = 0
$data[16,15] = 0
$data[16,16] = 0
$data[16,17] = 0.125
$data[17,0] = 0
$data[17,1] = 0
$data[17,2] = 0
$data[17,3] = 0
$data[17,4] = 0.01154529307282416
$data[17,5] = 0
$data[17,6] = 0.2007104795737123
$data[17,7] = 0.1012433392539965
$data[17,8] = 0.3641207815275311
$data[17,9] = 0.1287744227353464
$data[17,10] = 0
$data[17,11] = 0.01509769094138544
$data[17,12] = 0.0008880994671403197
$data[17,13] = 0.0630550621669627
$data[17,14] = 0
$data[17,15] = 0
$data[17,16] = 0
$data[17,17] = 0.1145648312611012

$ws.Range("B2:S19").Value2 = $data
